# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the
# 794b1975.../bdcc3160... row pair on both the zh-cn and de-de
# report sheets, reflecting a re-run of the handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-14 00:27:14"
$wsZhCn.Range("E5").Value = "2016-03-14 00:27:14"
$wsZhCn.Range("H4").Value = "2016-03-14 00:27:31"
$wsZhCn.Range("H5").Value = "2016-03-14 00:27:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-14 00:27:18"
$wsDeDe.Range("E5").Value = "2016-03-14 00:27:18"
$wsDeDe.Range("H4").Value = "2016-03-14 00:27:37"
$wsDeDe.Range("H5").Value = "2016-03-14 00:27:37"
